$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("material_granite_boulder", "Granite Boulder"),
    @("material_cotton_pillow", "Cotton Pillow"),
    @("material_rubber_duck", "Rubber Duck"),
    @("material_oak_branch", "Oak Branch")
)

$row = 16
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("B19").Select()
